# Data-driven parameter change across the Order Information (SMP) workbook.
# - Surface treatment dropdown value renamed: "Plasmanitrieren" -> "Plasma nitriding"
# - Material name renamed: "AlMg3-Blech" -> "EN AW-5754 / AlMg3"
# - Shipping option labels renamed:
#     "Standard shipping" -> "Package delivery (extra costs)"
#     "Pickup at factory" -> "Pick-up at factory (no costs)"
#     "Special packaging / via freight forwarding" -> "Freight delivery / sepcial packaging (extra costs)"
# - A couple of quantity values changed on the FPA012-013-015-017 sheet
# - Active sheet/selection moved from FPA011 to FPA012-013-015-017

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# FPA011
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("FPA011")
$ws.Range("E2").Value = "Plasma nitriding"
$ws.Range("M3").Value = "EN AW-5754 / AlMg3"
$ws.Range("P2").Value = "Package delivery (extra costs)"
$ws.Range("P3").Value = "Pick-up at factory (no costs)"
$ws.Range("P4").Value = "Freight delivery / sepcial packaging (extra costs)"
$ws.Range("P5").Value = "Package delivery (extra costs)"

# ---------------------------------------------------------------------------
# FPA012-013-015-017
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("FPA012-013-015-017")
$ws2.Range("D2").Value = 100
$ws2.Range("D3").Value = 200
$ws2.Range("D4").Value = 100
$ws2.Range("D5").Value = 100
$ws2.Range("E2").Value = "Plasma nitriding"
$ws2.Range("M3").Value = "EN AW-5754 / AlMg3"
$ws2.Range("P2").Value = "Package delivery (extra costs)"
$ws2.Range("P3").Value = "Pick-up at factory (no costs)"
$ws2.Range("P4").Value = "Freight delivery / sepcial packaging (extra costs)"
$ws2.Range("P5").Value = "Package delivery (extra costs)"

# ---------------------------------------------------------------------------
# FPA014-016-020
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("FPA014-016-020")
$ws3.Range("E2").Value = "Plasma nitriding"
$ws3.Range("P2").Value = "Package delivery (extra costs)"
$ws3.Range("P3").Value = "Pick-up at factory (no costs)"

# ---------------------------------------------------------------------------
# FPA018-019
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("FPA018-019")
$ws4.Range("E2").Value = "Plasma nitriding"
$ws4.Range("M3").Value = "EN AW-5754 / AlMg3"
$ws4.Range("Q2").Value = "Package delivery (extra costs)"
$ws4.Range("Q3").Value = "Pick-up at factory (no costs)"
$ws4.Range("Q4").Value = "Freight delivery / sepcial packaging (extra costs)"
$ws4.Range("Q5").Value = "Package delivery (extra costs)"

# ---------------------------------------------------------------------------
# BTMI010
# ---------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("BTMI010")
$ws5.Range("E2").Value = "Plasma nitriding"
$ws5.Range("M3").Value = "EN AW-5754 / AlMg3"

# ---------------------------------------------------------------------------
# BTMI016
# ---------------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item("BTMI016")
$ws6.Range("E2").Value = "Plasma nitriding"

# ---------------------------------------------------------------------------
# Move the active tab / selection to FPA012-013-015-017 (was FPA011)
# ---------------------------------------------------------------------------
$ws2.Activate()
$ws2.Range("D12").Select()
$ws.Range("J18").Select()
